$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (row 5) - the new dataset only has 3 data rows
$ws.Rows.Item(5).Delete()

# Row 2: FAPs / Cxcl13 / Ackr4 / ECs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Cxcl13"
$ws.Range("C2").Value = "Ackr4"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 10.91523966666667
$ws.Range("H2").Value = 32.745719
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.07058066666666667
$ws.Range("N2").Value = 0.211742
$ws.Range("O2").Value = 0.03405243646146196
$ws.Range("P2").Value = 0.03405243646146196
$ws.Range("Q2").Value = 0.7704048924997777
$ws.Range("R2").Value = 6.933644032498001
$ws.Range("S2").Value = 0.03405243646146196
$ws.Range("T2").Value = 0.03405243646146196

# Row 3: FAPs / Cxcl13 / Ackr4 / FAPs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Cxcl13"
$ws.Range("C3").Value = "Ackr4"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 10.91523966666667
$ws.Range("H3").Value = 32.745719
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 1.890486333333333
$ws.Range("N3").Value = 5.671459
$ws.Range("O3").Value = 0.9120863940138783
$ws.Range("P3").Value = 0.9120863940138783
$ws.Range("Q3").Value = 20.63511141489122
$ws.Range("R3").Value = 185.716002734021
$ws.Range("S3").Value = 0.9120863940138783
$ws.Range("T3").Value = 0.9120863940138783

# Row 4: FAPs / Cxcl13 / Ackr4 / MuSCs
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Cxcl13"
$ws.Range("C4").Value = "Ackr4"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 10.91523966666667
$ws.Range("H4").Value = 32.745719
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.1116383333333333
$ws.Range("N4").Value = 0.334915
$ws.Range("O4").Value = 0.05386116952465989
$ws.Range("P4").Value = 0.05386116952465989
$ws.Range("Q4").Value = 1.218559164320556
$ws.Range("R4").Value = 10.967032478885
$ws.Range("S4").Value = 0.05386116952465989
$ws.Range("T4").Value = 0.05386116952465989
